# Apply cryptos list update (GitHub Actions scheduled refresh)
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '65.336.97'
$ws.Range("E2").Value = '  +0.53%  '
$ws.Range("D3").Value = '3.212.40'
$ws.Range("E3").Value = '  -0.95%  '
$ws.Range("E4").Value = '  +0.00%  '
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '576.48'
$ws.Range("E5").Value = '  -0.44%  '
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '167.56'
$ws.Range("E6").Value = '  -3.39%  '
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = '0.598'
$ws.Range("E7").Value = '  -5.68%  '
$ws.Range("E8").Value = '  +0.08%  '
$ws.Range("E9").Value = '  -2.41%  '
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = '6.75'
$ws.Range("E10").Value = '  -0.61%  '
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = '0.391'
$ws.Range("E11").Value = '  +0.10%  '
$ws.Range("D12").Value = '3.771.58'
$ws.Range("E12").Value = '  -0.85%  '
$ws.Range("D14").Value = '65.319.71'
$ws.Range("E14").Value = '  +0.42%  '
$ws.Range("E15").Value = '  -0.63%  '
$ws.Range("D16").Value = '3.204.85'
$ws.Range("E16").Value = '  -0.83%  '
$ws.Range("E17").Value = '  -1.14%  '
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = '414.50'
$ws.Range("E18").Value = '  -1.08%  '
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = '12.94'
$ws.Range("E19").Value = '  +0.35%  '
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = '5.35'
$ws.Range("E20").Value = '  -1.00%  '
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = '7.18'
$ws.Range("E21").Value = '  -0.57%  '
$ws.Range("E22").Value = '  +0.10%  '
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = '69.78'
$ws.Range("E23").Value = '  -1.86%  '
$ws.Range("E24").Value = '  -1.41%  '
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = '0.491'
$ws.Range("E25").Value = '  -1.29%  '
$ws.Range("E26").Value = '  -5.69%  '
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = '8.94'
$ws.Range("E27").Value = '  -2.07%  '
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = '0.999'
$ws.Range("E28").Value = '  +0.02%  '
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = '1.84'
$ws.Range("E29").Value = '  -1.61%  '
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = '21.68'
$ws.Range("E30").Value = '  -1.09%  '
$ws.Range("E31").Value = '  -0.15%  '
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = '6.42'
$ws.Range("E32").Value = '  -0.63%  '
$ws.Range("E33").Value = '  -1.55%  '
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = '157.39'
$ws.Range("E34").Value = '  -0.01%  '
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = '1.38'
$ws.Range("E35").Value = '  -1.43%  '
$ws.Range("E36").Value = '  -0.21%  '
$ws.Range("D37").Value = '2.742.89'
$ws.Range("E37").Value = '  -3.38%  '
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = '24.26'
$ws.Range("E38").Value = '  -5.16%  '
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = '4.16'
$ws.Range("E39").Value = '  -2.18%  '
$ws.Range("E40").Value = '  -1.40%  '
$ws.Range("E41").Value = '  +0.48%  '
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = '5.61'
$ws.Range("E42").Value = '  -2.84%  '
$ws.Range("B43").Value = 'VeChain'
$ws.Range("C43").Value = 'https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet'
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = '0.0264'
$ws.Range("E43").Value = '  -0.57%  '
$ws.Range("B44").Value = 'Bittensor'
$ws.Range("C44").Value = 'https://coinranking.com/coin/pgv7xSFi6+bittensor-tao'
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = '297.17'
$ws.Range("E44").Value = '  -2.30%  '
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = '21.64'
$ws.Range("E45").Value = '  -2.61%  '
$ws.Range("E46").Value = '  +0.04%  '
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = '0.0992'
$ws.Range("E47").Value = '  -2.31%  '
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = '1.98'
$ws.Range("E48").Value = '  -8.78%  '
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = '5.82'
$ws.Range("E49").Value = '  -0.37%  '
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = '10.47'
$ws.Range("E50").Value = '  +0.61%  '
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = '0.910'
$ws.Range("E51").Value = '  -2.57%  '
